# --- Update "Data" sheet (sheet 1): refresh WTREGEN weekly series window ---
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New date serials and values (weekly WTREGEN observations) for rows 2..93
$dates = @(44475,44482,44489,44496,44503,44510,44517,44524,44531,44538,44545,44552,44559,44566,44573,44580,44587,44594,44601,44608,44615,44622,44629,44636,44643,44650,44657,44664,44671,44678,44685,44692,44699,44706,44713,44720,44727,44734,44741,44748,44755,44762,44769,44776,44783,44790,44797,44804,44811,44818,44825,44832,44839,44846,44853,44860,44867,44874,44881,44888,44895,44902,44909,44916,44923,44930,44937,44944,44951,44958,44965,44972,44979,44986,44993,45000,45007,45014,45021,45028,45035,45042,45049,45056,45063,45070,45077,45084,45091,45098,45105,45112)
$values = @(135.582,77.858,83.386,213.291,280.205,259.293,211.453,164.087,152.651,115.117,104.798,84.95399999999999,211.849,389.751,456.401,488.697,614.877,675.773,686.384,692.5410000000001,697.838,682.903,652.039,570.606,614.616,575.0650000000001,574.451,547.308,711.4059999999999,944.328,945.478,945.658,886.9640000000001,821.535,789.532,723.384,656.87,758.283,756.627,700.457,652.572,613.878,602.9450000000001,594.115,557.265,545.321,535.2670000000001,612.5359999999999,581.295,593.808,692.496,689.569,633.939,608.302,607.199,634.548,598.544,527.479,502.982,479.474,511.474,432.335,355.517,459.78,427.926,423.625,372.34,339.018,491.848,560.0890000000001,493.277,490.379,477.333,381.245,333.35,232.866,259.587,183.577,168.793,109.208,166.555,291.702,269.216,197.666,116.22,61.952,48.954,44.756,102.118,276.85,390.571,415.441)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Remove now-obsolete trailing rows (94..132) so the sheet dimension shrinks to A1:B93
$ws.Rows("94:132").Delete()

# --- Update "SeriesInfo" sheet (sheet 2): refresh FRED series metadata ---
$ws2 = $wb.Worksheets.Item("SeriesInfo")

# Use a leading quote so these date-shaped strings are kept as literal text,
# matching the original inline-string cells instead of being auto-parsed as dates.
$q = [char]39
$ws2.Range("B3").Value = $q + "2023-07-09"
$ws2.Range("B4").Value = $q + "2023-07-09"
$ws2.Range("B7").Value = $q + "2023-07-05"
$ws2.Range("B14").Value = "2023-07-06 15:34:05-05"
$ws2.Range("B15").Value = 84

